$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "sv"
$ws.Range("H12").Value = "vsfv"
$ws.Range("K7").Value = "vrwsf"

$ws.Range("K7").Select()
